$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (date style) from A8 down into A9:A10, matching the existing column style
$ws.Cells.Item(8, 1).Copy()
$ws.Range($ws.Cells.Item(9, 1), $ws.Cells.Item(10, 1)).PasteSpecial(-4122)

# New row 9 data
$ws.Cells.Item(9, 1).Value = 42612.883043981485
$ws.Cells.Item(9, 2).Value = 6
$ws.Cells.Item(9, 3).Value = 52
$ws.Cells.Item(9, 4).Value = 45
$ws.Cells.Item(9, 5).Value = 54
$ws.Cells.Item(9, 6).Value = 45
$ws.Cells.Item(9, 7).Value = 14346
$ws.Cells.Item(9, 8).Value = 13154
$ws.Cells.Item(9, 9).Value = 2379
$ws.Cells.Item(9, 10).Value = 260
$ws.Cells.Item(9, 11).Value = 226
$ws.Cells.Item(9, 12).Value = 6
$ws.Cells.Item(9, 13).Value = 5
$ws.Cells.Item(9, 14).Value = "Noun"

# New row 10 data
$ws.Cells.Item(10, 1).Value = 42612.88957175926
$ws.Cells.Item(10, 2).Value = 16
$ws.Cells.Item(10, 3).Value = 53
$ws.Cells.Item(10, 4).Value = 46
$ws.Cells.Item(10, 5).Value = 78
$ws.Cells.Item(10, 6).Value = 22
$ws.Cells.Item(10, 7).Value = 15158
$ws.Cells.Item(10, 8).Value = 14183
$ws.Cells.Item(10, 9).Value = 2635
$ws.Cells.Item(10, 10).Value = 297
$ws.Cells.Item(10, 11).Value = 256
$ws.Cells.Item(10, 12).Value = 39
$ws.Cells.Item(10, 13).Value = 11
$ws.Cells.Item(10, 14).Value = "Noun"
